$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as plain text so that
# numeric-looking values (e.g. "114.63") are not reinterpreted as numbers
# and values with thousand separators (e.g. "45.641.43") keep their exact
# string representation, matching the original inlineStr cell contents.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "45.641.43"
$ws.Range("E2").Value = "  +7.46%  "
$ws.Range("D3").Value = "2.386.00"
$ws.Range("E3").Value = "  +4.44%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "114.63"
$ws.Range("E5").Value = "  +11.16%  "
$ws.Range("D6").Value = "317.61"
$ws.Range("E6").Value = "  +2.32%  "
$ws.Range("D7").Value = "0.634"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.629"
$ws.Range("E9").Value = "  +4.62%  "
$ws.Range("D10").Value = "43.20"
$ws.Range("E10").Value = "  +11.66%  "
$ws.Range("E11").Value = "  +4.75%  "
$ws.Range("D12").Value = "8.71"
$ws.Range("E12").Value = "  +6.35%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("E14").Value = "  +4.66%  "
$ws.Range("E15").Value = "  +4.57%  "
$ws.Range("D16").Value = "2.746.74"
$ws.Range("E16").Value = "  +4.35%  "
$ws.Range("D17").Value = "2.388.09"
$ws.Range("E17").Value = "  +4.49%  "
$ws.Range("D18").Value = "45.556.51"
$ws.Range("E18").Value = "  +7.37%  "
$ws.Range("E19").Value = "  +4.20%  "
$ws.Range("E20").Value = "  +4.36%  "
$ws.Range("D21").Value = "13.33"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").Value = "74.84"
$ws.Range("E22").Value = "  +2.24%  "
$ws.Range("D23").Value = "3.52"
$ws.Range("E23").Value = "  +3.77%  "
$ws.Range("D24").Value = "269.19"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E25").Value = "  +9.68%  "
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "7.64"
$ws.Range("E27").Value = "  +10.16%  "
$ws.Range("E28").Value = "  +5.70%  "
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("D30").Value = "22.95"
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("D31").Value = "38.84"
$ws.Range("E31").Value = "  +8.60%  "
$ws.Range("D32").Value = "0.0968"
$ws.Range("E32").Value = "  +14.75%  "
$ws.Range("D33").Value = "171.19"
$ws.Range("E33").Value = "  +4.26%  "
$ws.Range("E34").Value = "  +17.55%  "
$ws.Range("E35").Value = "  +11.21%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.120"
$ws.Range("E36").Value = "  +8.28%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "0.132"
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("D38").Value = "4.18"
$ws.Range("E38").Value = "  +15.84%  "
$ws.Range("D39").Value = "3.05"
$ws.Range("E39").Value = "  +11.44%  "
$ws.Range("D40").Value = "0.0366"
$ws.Range("E40").Value = "  +6.34%  "
$ws.Range("D41").Value = "1.72"
$ws.Range("E41").Value = "  +11.36%  "
$ws.Range("D42").Value = "103.61"
$ws.Range("E42").Value = "  -7.31%  "
$ws.Range("D43").Value = "0.239"
$ws.Range("E43").Value = "  +6.81%  "
$ws.Range("D44").Value = "71.47"
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("D45").Value = "13.28"
$ws.Range("E45").Value = "  +10.58%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "5.79"
$ws.Range("E47").Value = "  +13.05%  "
$ws.Range("D48").Value = "116.26"
$ws.Range("E48").Value = "  +5.81%  "
$ws.Range("E49").Value = "  +17.06%  "
$ws.Range("D50").Value = "9.38"
$ws.Range("E50").Value = "  +8.41%  "
$ws.Range("D51").Value = "79.07"
$ws.Range("E51").Value = "  +2.74%  "
